$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E in this sheet are stored as text (inline strings) even
# though many values look numeric (e.g. "92.43", "302.30", "41.987.75").
# Force text format on D-column cells before assigning so Excel keeps the
# literal text (preserving trailing zeros / thousands-style dots) instead of
# re-interpreting the input as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.023.31"
$ws.Range("E2").Value = "  +5.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.259.92"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.30"
$ws.Range("E5").Value = "  +3.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.43"
$ws.Range("E6").Value = "  +6.00%  "
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.68"
$ws.Range("E10").Value = "  +9.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.36"
$ws.Range("E11").Value = "  +6.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0800"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.609.61"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.09"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.267.05"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.904.45"
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.17"
$ws.Range("E20").Value = "  +9.51%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.02"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.06"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.90"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.96"
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("E29").Value = "  +13.22%  "
$ws.Range("E30").Value = "  +4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.45"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.94"
$ws.Range("E32").Value = "  +6.12%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0745"
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("E38").Value = "  +3.68%  "
$ws.Range("E39").Value = "  +4.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.51"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("E42").Value = "  +5.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.057.33"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.74"
$ws.Range("E44").Value = "  +9.35%  "
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.87"
$ws.Range("E47").Value = "  +5.26%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.61"
$ws.Range("E51").Value = "  +5.14%  "
